$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.079.10"
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").Value = "1.551.10"
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").Value = "'287.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("D7").Value = "'0.3826"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.60%  "

$ws.Range("D8").Value = "'0.3276"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.41%  "

$ws.Range("D9").Value = "'43.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.64%  "

$ws.Range("E10").Value = "  -0.75%  "

$ws.Range("D11").Value = "'0.07351"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "'20.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.92%  "

$ws.Range("D14").Value = "'5.784"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.60%  "

$ws.Range("D15").Value = "'6.752"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.35%  "

$ws.Range("D16").Value = "1.562.83"
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("D17").Value = "'0.00001083"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.83%  "

$ws.Range("D18").Value = "'0.06633"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.75%  "

$ws.Range("D19").Value = "'85.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.38%  "

$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").Value = "'6.349"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "'16.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.47%  "

$ws.Range("D23").Value = "'11.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.29%  "

$ws.Range("D24").Value = "22.070.80"
$ws.Range("E24").Value = "  -1.71%  "

$ws.Range("D25").Value = "'2.295"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.29%  "

$ws.Range("D26").Value = "'2.499"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.63%  "

$ws.Range("D27").Value = "'150.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.57%  "

$ws.Range("E28").Value = "  -2.89%  "

$ws.Range("D29").Value = "'4.930"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.81%  "

$ws.Range("D30").Value = "'121.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.00%  "

$ws.Range("D31").Value = "1.736.53"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("D32").Value = "'1.078"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.17%  "

$ws.Range("D33").Value = "'5.865"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.22%  "

$ws.Range("D34").Value = "'1.900"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.42%  "

$ws.Range("D35").Value = "'0.08221"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.90%  "

$ws.Range("D36").Value = "'9.257"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.40%  "

$ws.Range("D37").Value = "'0.06289"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.66%  "

$ws.Range("D38").Value = "'0.02312"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.92%  "

$ws.Range("D39").Value = "'5.259"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.05%  "

$ws.Range("E40").Value = "  -5.75%  "

$ws.Range("D41").Value = "'1.229"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.53%  "

$ws.Range("D42").Value = "'11.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.33%  "

$ws.Range("D43").Value = "'1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").Value = "'0.6002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.50%  "

$ws.Range("D45").Value = "'13.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.29%  "

$ws.Range("E46").Value = "  -1.16%  "

$ws.Range("D47").Value = "'0.5813"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.18%  "

$ws.Range("D48").Value = "'1.966"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.84%  "

$ws.Range("D49").Value = "'121.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.91%  "

$ws.Range("E50").Value = "  -3.19%  "

$ws.Range("D51").Value = "'0.07011"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.01%  "
